# Petty cash book update — 20-Mei-2021, end of day update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 18: BENSIN top-up added to existing Wages Expense formula ---
$ws.Range("D18").Formula = "=45000+195000"

# --- Row 19: add a further TRANSFER BCA component ---
$ws.Range("D19").Formula = "=3640000+1250000+1170000+440000+515000+2140000"

# --- Row 20: add further A/R receipt component ---
$ws.Range("C20").Formula = "=9300000+17306000"

# --- New transactions appended for 20-Mei-2021 ---

# Row 23: BENSIN - RUSH
$ws.Range("B23").Value = "BENSIN - RUSH"
$ws.Range("D23").Formula = "=250000"

# Row 24: JASON - paspor
$ws.Range("B24").Value = "JASON - paspor"
$ws.Range("D24").Value = 1350000

# Row 25: SALES - cash/retail
$ws.Range("B25").Value = "SALES - cash/retail"
$ws.Range("C25").Formula = "=10414475+15258525-17306000"

# Row 26: SELISIH - lebih
$ws.Range("B26").Value = "SELISIH - lebih"
$ws.Range("C26").Value = 2000

# Row 27: SETOR KE BANK
$ws.Range("B27").Value = "SETOR KE BANK"
$ws.Range("D27").Formula = "=15000000"

# Row 28: new day, 20-May-2021 (serial 44336), Wages Expense
$ws.Range("A28").Value = 44336
$ws.Range("B28").Value = "Wages Expense"
$ws.Range("D28").Formula = "=45000"

# Row 29: A/R
$ws.Range("B29").Value = "A/R"
$ws.Range("C29").Formula = "=3360000"

# Row 30: TRANSFER BCA
$ws.Range("B30").Value = "TRANSFER BCA"
$ws.Range("D30").Formula = "=3360000"

# Row 31: A/P
$ws.Range("B31").Value = "A/P"
$ws.Range("D31").Formula = "=1877000"

# Row 32: FREIGHT OUT
$ws.Range("B32").Value = "FREIGHT OUT"
$ws.Range("D32").Formula = "=14500"

# --- Restore view state: scroll so row 27 is visible under the frozen pane, ---
# --- and leave the active selection on C48 (matches the author's end-of-day position). ---
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C48").Select()

$wb.Application.Calculate()
